# ---------------------------------------------------------------------------
# Applies the commit "feat: add 2022-Q3 data":
#   1. Inserts a new summary row (row 2) in the "总计" sheet for 2022-Q3
#      (10 -> 31 holdings, 6.07 -> 7.13 yi), shifting the prior rows down
#      and re-numbering the A-column sequence index.
#   2. Inserts a brand-new worksheet named "2022-Q3" right after "总计",
#      holding the per-fund holdings detail for the new quarter (31 funds).
# All other existing quarter sheets (2022-Q2 .. 2020-Q4) keep their names
# and data; they simply shift one tab to the right because of the insert.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1) "总计" (summary) sheet: insert the new 2022-Q3 row at row 2
# ===========================================================================
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
# Excel's row-insert inherits the row-1 header formatting onto the new
# row for the columns that had none before; strip that back off so B2:D2
# come out unstyled, matching the rest of the data rows.
$summary.Range("B2:D2").ClearFormats()

# Column A carries the bordered/bold "index" style (same as every other
# row in this column) - clone it from the row below instead of hand
# rebuilding the style.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 31
$summary.Range("D2").Value = 7.13

# Renumber the sequence index (column A) of the rows that got pushed down.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ===========================================================================
# 2) New "2022-Q3" worksheet, positioned right after "总计"
# ===========================================================================
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund holdings detail rows: idx, code, name, size, position, pct, value, rank
$fundData = @(
    ,@('0','000979','景顺长城沪港深精选股票','20.32','80.04','8.01','1.6276','2')
    ,@('1','260112','景顺长城能源基建混合A','22.22','61.72','6.69','1.4865','5')
    ,@('2','008850','景顺长城价值稳进三年定期开放灵活配置混合','17.29','84.84','5.71','0.9873','7')
    ,@('3','009098','景顺长城价值领航两年持有期混合','7.16','90.92','8.98','0.6430','5')
    ,@('4','008715','景顺长城价值驱动一年持有期灵活配置混合','8.44','90.91','7.08','0.5976','7')
    ,@('5','008060','景顺长城价值边际灵活配置混合A','5.45','80.34','7.98','0.4349','3')
    ,@('6','161611','融通内需驱动混合A/B','8.72','90.68','4.31','0.3758','6')
    ,@('7','014109','融通内需驱动混合C','4.06','90.68','4.31','0.1750','6')
    ,@('8','001319','农银汇理信息传媒主题股票','2.62','80.55','4.31','0.1129','8')
    ,@('9','005561','创金合信中证红利低波动指数A','3.32','94.12','2.11','0.0701','7')
    ,@('10','002252','融通成长30灵活配置混合A/B','1.56','93.13','4.05','0.0632','10')
    ,@('11','015779','景顺长城价值边际灵活配置混合C','0.79','80.34','7.98','0.0630','3')
    ,@('12','012708','东方红中证东方红红利低波动指数A','3.27','93.80','1.78','0.0582','3')
    ,@('13','512890','华泰柏瑞中证红利低波动ETF','2.60','99.50','2.24','0.0582','7')
    ,@('14','014106','融通成长30灵活配置混合C','1.34','93.13','4.05','0.0543','10')
    ,@('15','005562','创金合信中证红利低波动指数C','2.19','94.12','2.11','0.0462','7')
    ,@('16','008115','天弘中证红利低波动100指数C','2.44','94.56','1.76','0.0429','7')
    ,@('17','009726','招商中证500等权重指数增强A','2.67','90.23','1.44','0.0384','4')
    ,@('18','008114','天弘中证红利低波动100指数A','1.89','94.56','1.76','0.0333','7')
    ,@('19','001223','鹏华文化传媒娱乐股票','0.76','86.75','4.24','0.0322','6')
    ,@('20','515100','景顺长城中证红利低波动100ETF','1.62','98.63','1.84','0.0298','7')
    ,@('21','009658','汇丰晋信中小盘低波动策略股票A','0.85','90.14','2.02','0.0172','1')
    ,@('22','009727','招商中证500等权重指数增强C','1.12','90.23','1.44','0.0161','4')
    ,@('23','012709','东方红中证东方红红利低波动指数C','0.67','93.80','1.78','0.0119','3')
    ,@('24','007751','景顺长城中证沪港深红利成长低波动指数A','0.67','90.27','1.71','0.0115','4')
    ,@('25','003359','大成中证360互联网+大数据100指数C','1.11','92.17','0.99','0.0110','6')
    ,@('26','011824','浙商汇金量化臻选股票A','0.88','92.26','1.20','0.0106','9')
    ,@('27','002236','大成中证360互联网+大数据100指数A','1.03','92.17','0.99','0.0102','6')
    ,@('28','011825','浙商汇金量化臻选股票C','0.39','92.26','1.20','0.0047','9')
    ,@('29','007760','景顺长城中证沪港深红利成长低波动指数C','0.06','90.27','1.71','0.0010','4')
    ,@('30','009775','汇丰晋信中小盘低波动策略股票C','0.04','90.14','2.02','0.0008','1')
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $i + 2
    $rec = $fundData[$i]

    $q3.Cells.Item($row, 1).Value = [int]$rec[0]

    # Columns B-G are stored as *text* in the source data (fund code,
    # name, size, position, pct, market value) even though several look
    # numeric - force text storage via NumberFormat before assignment so
    # Excel doesn't silently coerce them to numbers.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q3.Cells.Item($row, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$c - 1]
    }

    $q3.Cells.Item($row, 8).Value = [int]$rec[7]
}

# Header row (B1:H1) and the index column (A2:A32) use the same
# bold/bordered/centered style used throughout the workbook - clone it
# from the "总计" sheet instead of re-deriving the raw style index.
$summary.Range("A2").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A32").PasteSpecial(-4122)

# Re-apply the header text and indices (PasteSpecial(-4122) only touches
# formats, so the values written above should be untouched, but re-set
# just in case of any clipboard interaction).
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $i + 2
    $rec = $fundData[$i]
    $q3.Cells.Item($row, 1).Value = [int]$rec[0]
    $q3.Cells.Item($row, 8).Value = [int]$rec[7]
}

$q3.Range("A1").Select()
